# Update latest output (run 129)
# Refresh the "Schedule" and "Detailed" sheets with the newest
# optimisation-run values/statuses.

$wb = $excel.ActiveWorkbook

$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsSchedule.Range("B4").Value = 46043.125
$wsSchedule.Range("C4").Value = 5
$wsSchedule.Range("D4").Value = 18.9
$wsSchedule.Range("E4").Value = 562.6022609999999
$wsSchedule.Range("F4").Value = 29.76731539682539
$wsSchedule.Range("A5").Value = 46043.29166666666
$wsSchedule.Range("C5").Value = 9
$wsSchedule.Range("D5").Value = 34.02
$wsSchedule.Range("E5").Value = -177.4568445
$wsSchedule.Range("F5").Value = -5.21625057319224

$wsDetailed = $wb.Worksheets.Item("Detailed")
$wsDetailed.Range("B35").Value = -5.76787
$wsDetailed.Range("B36").Value = -2.7996
$wsDetailed.Range("B37").Value = 49.38629
$wsDetailed.Range("B38").Value = 60.01523
$wsDetailed.Range("C38").Value = "historical"
$wsDetailed.Range("B39").Value = 63.25126
$wsDetailed.Range("B40").Value = 83.7666
$wsDetailed.Range("B41").Value = 81.25905
$wsDetailed.Range("B42").Value = 79.08551
$wsDetailed.Range("B43").Value = 78
$wsDetailed.Range("B44").Value = 75.13731
$wsDetailed.Range("B45").Value = 64.8901
$wsDetailed.Range("B47").Value = 63.00775
$wsDetailed.Range("E56").Value = "OFF"
$wsDetailed.Range("B58").Value = 65
$wsDetailed.Range("B59").Value = 68.17046999999999
$wsDetailed.Range("B60").Value = 73.20005
$wsDetailed.Range("B61").Value = 78.85684000000001
$wsDetailed.Range("B62").Value = 79.26575
$wsDetailed.Range("B63").Value = 64.89
$wsDetailed.Range("B64").Value = 54.86055
$wsDetailed.Range("E64").Value = "ON"
$wsDetailed.Range("B66").Value = 0.00002
$wsDetailed.Range("B67").Value = -5.50985
$wsDetailed.Range("B68").Value = -6.97897
$wsDetailed.Range("B69").Value = -7.1246
$wsDetailed.Range("B70").Value = -8.87745
$wsDetailed.Range("B71").Value = -9.5
$wsDetailed.Range("B73").Value = -14.76574
$wsDetailed.Range("B74").Value = -16.47786
$wsDetailed.Range("B75").Value = -22.30467
$wsDetailed.Range("B76").Value = -22.25351
$wsDetailed.Range("B77").Value = -24.46768
$wsDetailed.Range("B78").Value = -24.10291
$wsDetailed.Range("B79").Value = -22.49172
$wsDetailed.Range("B81").Value = -15.21263
$wsDetailed.Range("B82").Value = -6.71362
$wsDetailed.Range("B83").Value = -5.13153
$wsDetailed.Range("B85").Value = 47.13993
$wsDetailed.Range("B86").Value = 56.52442
$wsDetailed.Range("B87").Value = 64.8901
$wsDetailed.Range("B88").Value = 100.01
$wsDetailed.Range("B89").Value = 100.01
$wsDetailed.Range("B90").Value = 86.53725
$wsDetailed.Range("B91").Value = 73.19
$wsDetailed.Range("B92").Value = 65
$wsDetailed.Range("B94").Value = 59.91511
$wsDetailed.Range("B95").Value = 63.83478
